# Auto-generated edit script to set full data grid per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('RECALL_CLASSIFICATION_DATE', 'PRODUCT_TYPE', 'CLASSIFICATION', 'RECALL_NUMBER', 'PRODUCT', 'RECALLING_FIRM', 'MANUFACTURER', 'RECALL_INITIALIZATION_DATE', 'REASON', 'VOLUME', 'DISTRIBUTION'),
    @('1995-04-12', 'Drugs', 'III', 'd-106-5', 'A B Antipyrine Benzocaine Otic Drops Rx In 10 Ml And 1 2 Fluid Ounce Glass Bottles Under The Following Labels Goldline H', '5 Manufacturer', 'Clay Park Laboratories Inc', 'April 5 1995', 'Product Does Not Meet Ph And Water Content Specifications', 'Approximately 534 183 10 Ml Bottles And 800 797 15 Ml Glass Bottles Were Distributed', 'Nationwide Virgin Islands West Indies Guyana'),
    @('1995-04-12', 'Drugs', 'III', 'd- 107-5', 'Potassium Phosphates Injection Usp In 5 Ml And 50 Ml Single Dose Vials Used As A Source Of Phosphorus Recall D 107 5', 'Manufacturer', 'Luitpold Pharmaceuticals Inc', 'March 16 1995', 'Presence Of Crystallization', '115 325 Vials Of Lot 4399 And 26 450 Vials Of Lot 4531', 'Nationwide'),
    @('1995-04-12', 'Drugs', 'III', 'd-108-5', 'Rectacaine And All Other Private Labels Of The Firm''s Hemorrhoidal Suppositories Containing Live Yeast Cell Derivatives Lycd', 'The Reese Chemical Company', 'Dr', 'March 28 1995', 'Product Distributed After September 1994 Ruling Regarding Live Yeast Cell Derivative Lycd Content', '323 Packages Containing 12 Suppositories Each Were Distributed', 'Nationwide'),
    @('1995-04-12', 'Drugs', 'III', 'd-109-5', 'Bisacodyl Suppositories 10 Mg Usp', 'Manufacturer', 'Clay Park Laboratories Inc', 'April 17 1995', 'Subpotenc', 'A Total Of 758 339 Labeled Packages Were Distributed To Date 427 466 8 S 272 014 12 S 7 526 50 S 49 986 100 S 1 171 500 S 176 1000 S This Product Was Also Sold In Bulk To Two Customers E Z Em Inc And P F Laboratories', 'Nationwide'),
    @('1995-04-12', 'Drugs', 'III', 'd-110-5', 'Supac Analgesic Tablets In 100 And 1000 Count Bottles Otc Pain Reliever', 'Manufacturer', 'Mission Pharmacal Company San Antonio Texas', 'March 23 1995', 'Product Does Not Meet Uniformity Specifications For The Caffeine Ingredient', '2 964 100 Tablet And 299 1000 Tablet Bottles Were Distributed Firm Estimated That 20 30 Percent Of Product Remained On Market At Time Of Recall', 'Nationwide Sri Lanka'),
    @('1995-04-12', 'Drugs', 'III', 'd-111-5', 'Hemorrhoidal A Suppositories 2', 'Manufacturer', 'Clay Park Laboratories Inc', 'March 30 1995', 'Stability Data Does Not Support Expiration Date', 'Approximately 75 185 Units Were Distributed', 'New York New Jersey California Ohio Florida Georgia Texas Pennsylvania Mississippi Tennessee Massachusetts Arkansas Alabama Indiana Ohio Puerto Rico Hong Kong'),
    @('1995-04-12', 'Drugs', 'III', 'b-198-5', 'Human Tissue From Cadaveric Donors Indicated For Human Transplant', 'Manufacturer', 'Lions Doheny Eye And Tissue Bank Los Angeles California', 'November 1 1994', 'Human Tissue From Cadaveric Donors Was Distributed Which 1 Lacked Adequate Documentation About The Donors Relevant Medical History Related To Risk Factors For Or Clinical Evidence Of Hepatitis B Hepatitis C Or Human Immunodeficiency Virus Hiv Infection Or 2 The Donor Tested Repeatedly Reactive For Hepatitis B Surface Antigen Hbsag', 'All Lots Involved', 'California Florida Indiana New Jersey New York Georgia Washington D'),
    @('1995-04-12', 'Biologics', 'II', 'b-195/196-5', 'A Red Blood Cells B Platelets', 'Manufacturer', 'Central Florida Blood Bank Inc', 'July 1 1994', 'Blood Products Which Were Collected From A Donor Who Had A History Of Nose Piercing Within 12 Months Of Donation Were Distributed', '1 Unit Of Each Component', 'Florida'),
    @('1995-04-12', 'Biologics', 'II', 'b-201/202-5', 'A Red Blood Cells B Recovered Plasma', 'Manufacturer', 'Manatee Community Blood Center Inc', 'November 29 1994', 'Blood Products Which Tested Initially Reactive For The Hepatitis B Surface Antigen Hbsag Followed By A Single Repeat Negative Test For Hbsag Were Distributed', '1 Unit Of Each Component', 'Florida'),
    @('1995-04-12', 'Biologics', 'II', 'b-203-5', 'Fresh Frozen Plasma', 'Manufacturer', 'South Bend Medical Foundation Inc', 'August 4 1992', 'Blood Product Which Tested Repeatedly Reactive For The Antibody To The Human Immunodeficiency Virus Type 1 And 2 Anti Hiv 1 2 Was Distributed', '1 Unit', 'Indiana'),
    @('1995-04-12', 'Biologics', 'II', 'b-204-5', 'Whole Blood Cpda 1', 'Manufacturer', 'Community Blood Center Of Greater Kansas City Kansas City Missouri', 'May 25 1994', 'Blood Product Which Tested Negative For The Antibody To The Human Immunodeficiency Virus Type 1 Anti Hiv 1 But Was Collected From A Donor Who Previously Tested Repeatedly Reactive To Anti Hiv 1 Was Distributed', '1 Unit', 'Missouri'),
    @('1995-04-12', 'Biologics', 'II', 'b-206-5', 'Platelets Pheresis', 'Manufacturer', 'South Bend Medical Foundation Inc', 'June 29 1992', 'Blood Product Which Tested Negative For The Antibody To The Human Immunodeficiency Virus Type 1 And 2 Anti Hiv 1 2 But Was Collected From A Donor Who Previously Tested Repeatedly Reactive For Anti Hiv 1 2 Was Distributed', '1 Unit', 'Indiana'),
    @('1995-04-12', 'Biologics', 'II', 'b-207-5', 'Red Blood Cells', 'Manufacturer', 'Walter L', 'November 7 1994', 'Blood Product Which Tested Negative For The Antibody To The Hepatitis B Core Antigen Anti Hbc But Was Collected From A Donor Who Previously Tested Repeatedly Reactive For Anti Hbc Was Distributed', '1 Unit', 'Georgia'),
    @('1995-04-12', 'Biologics', 'II', 'b-208/209-5', 'A Red Blood Cells B Platelets', 'Manufacturer', 'Walter L', 'February 10 1994', 'Blood Products Which Tested Initially Reactive For The Antibody To The Hepatitis B Core Antigen Anti Hbc Were Distributed', '1 Unit Of Each Component', 'Georgia'),
    @('1995-04-12', 'Biologics', 'II', 'b-218-5', 'Red Blood Cells', 'Manufacturer', 'Mississippi Valley Regional Blood Center Davenport Iowa', 'August 25 1994', 'Blood Product Which Was Improperly Tested For The Antibody To The Human Immunodeficiency Virus Type 1 Anti Hiv 1 Was Distributed', '1 Unit', 'Illinois'),
    @('1995-04-12', 'Devices', 'II', 'z-335/338-5', 'Siemens Mevasim Radiation Therapy Simulator Systems With Plastic Rotatable Cassette Holders Teletherapy Planning Systems For The Treatment Of Cancer A Mevasim Radiation Therapy Simulator System B Mevasim 2 Radiation Therapy Simulator System C Mevasim S Radiation Therapy Simulator System D Mevasim S 27 Radiation Therapy Simulator System', 'Siemens Medical Systems', 'Mecaserto St', 'May 23 1994', 'Assemblies Used To Attach The Rotatable Cassette Holder To A Mevasim System May Fail Causing The Cassette Holder Assembly To Fall Product 5 In 1 Sterile Plastic Tubing Connectors Single Use Recall Z 601 5 Code Lot Numbers 9776 9777 9778 9779 Representing All Lots Manufacturer Mcp Industries Inc Decatur Alabama Recalled By Sterile Design Division Of Maxxim Medical Clearwater Florida By Letter February 14 1995 Firm Initiated Recall Complete Distribution Alabama Arkansas California Florida Georgia Massachusetts Mississippi Missouri New Jersey North Carolina Pennsylvania South Carolina Tennessee Texas West Virginia Quantity Approximately 216 000 Devices Were Manufactured And Are Subject To Recall Reason The Sterility Of The Connectors Has Been Compromised Due To Inadequate Documentation To Assure Adequate Sterilization Product Mp Video Medicam 900 Insufflator Sterile Disposable Filter And Tubing Set Recall Z 602 5 Code Product 01 5212 Lot Af 4 L 0002 Manufacturer Cybermedic Boulder Colorado Recalled By Mp Video Medway Massachusetts By Telephone Or By Fax January 26 1995 Followed By Letter January 27 1995 Firm Initiated Recall Ongoing Distribution Connecticut Florida Illinois Massachusetts Maine North Carolina New Hampshire New York Texas Utah Wisconsin Chile Spain Colombia Belgium Quantity 63 Cases 10 Bags Per Case Were Distributed Reason The Sterility Of The Device Has Been Compromised As Evidenced By Loss Of Package Integrity Product Sterile Design Custom Surgical Procedure Trays Containing Coated Vicryl Sutures Note Only The Referenced Sutures Made By Ethicon Are Under Recall Not The Entire Surgical Tray Recall Z 603 635 5 Mod Cat Code Akron City Hosp 87 Ach 01 Ao 88 Ach 01 Ao 89 Ach 01 Ao 12 Baptist Hospital Nashville 70 Bh 33 Nt 68 Bh 33 Nt 71 Bh 33 Nt Baptist Med Ctr 48 Bm 04 Kcm 46 Bm 04 Kcm 47 Bm 04 Kcm 33 Bmc 12 Oco 31 Bmc 12 Oco 32 Bmc 12 Oco 35 Bmc 12 Oco 18 Bmc 23 Oco 19 Bmc 23 Oco 34 Bmc 12 Oco Community Hosp 03 Ch 07 Mi 03 Ch 09 Mi 03 Ch 08 Mi 04 Ch 08 Mi 03 Ch 09 Mi 04 Ch 07 Mi 04 Ch 09 Mi 05 Ch 08 Mi Cuyahoga Falls Hosp 03 Cfgh 17 Cfo 04 Cfgh 17 Cfo 05 Cfgh 15 Cfo 07 Cfgh 10 Cfo 08 Cfgh 10 Cfo 03 Cfgh 18 Cfo 06 Cfgh 12 Cfo 04 Cfgh 15 Cfo 02 Cfgh 18 Cfo Cuyahoga Falls Hosp 04 Cfgh 18 Cfo 01 Cfgh 19 Cfo Dr Michael Katz 17 Mk 01 Vf Elyria Memorial 02 Emh 03 Ec Geisinger Med Ctr 30 Gmc 09 Dp 57 Gmc 04 Dp 21 Gmc 10 Dp 22 Gmc 10 Dp 40 Gmc 06 Dp Gillette Childrens 18 Gch 06 Spm Heb Harris Med 94 Heb 01 Bt Heb Harris Meth 19 Heb 18 Bt 20 Heb 18 Bt 92 Heb 01 Bt 93 Heb 01 Bt Hillcrest Hospital 14 Hh 12 Mho 15 Hh 12 Mho Iowa Methodist 02 Imm 02 Di 03 Imm 02 Di Jennie Edmundson 26 Jeh 02 Cbi Jewish Hospital 12 Jh 01 Lk 13 Jh 01 Lk 14 Jh 01 Lk Lutheran General 01 Lgh 25 Pri 02 Lgh 25 Pri 03 Lgh 25 Pri Mendocino Coast 32 Mch 01 Fbc 33 Mch 01 Fbc Menorah Medical Ctr 69 Mmc 10 Kcm 50 Mmc 12 Kcm 51 Mmc 12 Kcm 70 Mmc 10 Kcm 49 Mmc 12 Kcm 68 Mmc 10 Kcm Methodist Southlake 52 Mh 03 Mi 53 Mh 03 Mi 86 Mh 10 Mi 37 Mh 13 Mi Mississippi Meth 25 Mmrc 05 Jm 24 Mmrc 07 Jm 26 Mmrc 05 Jm Missouri Baptist 94 Mbho 4 Slm 95 Mbh 04 Slm 96 Mbh 04 Slm N C Baptist Hosp 38 Ncbh 03 Wsnc 86 Ncbh 01 Wsnc 87 Ncbh 01 Wsnc 39 Ncbh 03 Wsnc 88 Ncbh 01 Wsnc 40 Ncbh 01 Wsnc 90 Ncbh 01 Wsnc Nix Med Ctr 14 Nmc 06 Sat 15 Nmc 06 Sat Sentara Leigh Hosp 71 Lmh 04 Nv St Anthony Med Ctr 49 Samc 11 Slm St Catherine''s Hosp 62 Sch 13 Eci St Clare Hosp 01 Sch 09 Tw St Joseph Hosp Mi 04 Sjm 04 Mi St Joseph Hosp 01 Sjh 11 Tf 13 St Mary Med Ctr 57 Smmc 05 Hi 58 Smmc 05 Hi 11 Smmc 06 Hi 56 Smmc 05 Hi St Marys Hosp 18 Smh 20 Rv 19 Smh 20 Rv 22 Smh 18 Rv 20 Smh 20 Rv 07 Smh 19 Rv Trinity Luth Hosp 78 Tlh 01 Kcm 79 Tlh 01 Kcm Williamsport Hosp 13 Whp 03 Wp Manufacturer Sterile Design Clearwater Florida Recalled By Manufacturer Letters Of September 16 And 19 1994 Firm Initiated Recall Complete Distribution Nationwide Quantity 1 873 Trays Were Distributed Reason The Sterility Of The Sutures Contained In These Devices Cannot Be Assured Because Of Possible Microorganism Contamination During Aseptic Processing Following Sterilization Product Iv Extension Set With Add On Filter Indicated For The Iv Infusion Of Medical Fluids Recall Z 636 5 Code Catalog V 5812 Lot F 4 K 307 Manufacturer Mcgaw Inc Sabana Grande Puerto Rico Recalled By Mcgaw Inc Irvine California By Letter December 22 1994 Firm Initiated Recall Ongoing Distribution Nationwide Canada Quantity 643 Cases Were Distributed Reason Due To An Inadequate Solvent Bond The Product Leaks Between The Tubing And Filter Assembly Which Could Lead To Under Infusion And Blood Loss Product Groshong Catheter Trays A 9 5 French Dual Lumen Groshong Catheter Tray B 7 French Single Lumen Groshong Catheter Tray Recall Z 640 641 5 Code Product No Lot No A 7726952 36 Ie 1909 B 60366 4 36 Ie 1909 Manufacturer Bard Access Systems Inc Salt Lake City Utah Recalled By Manufacturer By Telephone November 22 1994 Followed By Letter November 28 1994 Firm Initiated Recall Ongoing Distribution Alabama Arizona California Colorado Florida Illinois Michigan New Hampshire New Jersey Oklahoma Texas West Virginia Quantity 64 Units Were Distributed Reason Mislabeling Single Lumen Catheters Were Inadvertently Mixed In With Dual Lumen Catheters Update Recall Z 440 441 5 Zinetics 24 Single And Dual Channel Internal Reference Catheters Manufactured By Zinetics Medical Inc Salt Lake City Utah Which Appeared In The March 8 1995 Enforcement Report Has Been Extended To 14 Include Additional Lot Numbers Z 440 5 Part 15 245 Lot Z 910 Z 441 5 Part 15 247 15 Lot Numbers Z 917 Z 963 Update Recall Z 471 5 Pegasys 2 1 Software System Manufactured By Adac Laboratories Milpitas California Which Appeared In The March 15 1995 Enforcement Report Should Read Reason An Error In The Software Could Cause The Upper Portion Of The Image To Appear Brighter Than It Should Be The Problem Could Lead To A Misdiagnosis', '259 Units Are Subject To Field Correction', 'Nationwide And International'),
    @('1995-04-12', 'Devices', 'III', 'z-638-5', 'Rota Stat In Vitro Diagnostic Test Kits Used To Detect Rotavirus''s', 'Isolab', 'Eldan Technologies Company Ltd', 'February 1 1995', 'The Control Latex In This Kit Is Clumping Or Agglutinating', '41 Kits Were Distributed', 'Nationwide'),
    @('1995-04-12', 'Devices', 'III', 'z-639-5', 'Sportrol Growth Promotion Test Suspension Candida Albicans Used To Test The Growth Properties Of Prepared Bacteriological Media', 'Manufacturer', 'Namsa Northwood Ohio', 'January 30 1995', 'The Growth Promotion Test Suspension Is Contaminated With Bacillus Species', '35 Vials Were Distributed', 'Nationwide And Japan'),
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}
